# Fill in the previously-empty "idPagamento" (D column) values for rows 10-15.
# These values look numeric but must stay text cells (same as the rest of the
# column), so we briefly mark the cell as Text, assign the value, then put the
# cell style back to Normal so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D10") "76941312538"
Set-TextValue $ws.Range("D11") "76705028247"
Set-TextValue $ws.Range("D12") "76946688512"
Set-TextValue $ws.Range("D13") "76946729412"
Set-TextValue $ws.Range("D14") "76946761994"
Set-TextValue $ws.Range("D15") "76947984334"

# Append new rows 16-28, matching the shape of the existing "Vitor Ito" rows.
$newRows = @(
    @{ Row=16; C="11966548087"; D="76947869106" },
    @{ Row=17; C="11966548087"; D="76947967362" },
    @{ Row=18; C="11966548087"; D="76948023750" },
    @{ Row=19; C="11966548087"; D="76707382593" },
    @{ Row=20; C="11988776655"; D="76707312891" },
    @{ Row=21; C="11966548088"; D="76707460717" },
    @{ Row=22; C="11977665544"; D="76950373414" },
    @{ Row=23; C="11966548087"; D="" },
    @{ Row=24; C="11988776655"; D="" },
    @{ Row=25; C="11977665544"; D="" },
    @{ Row=26; C="11977665544"; D="" },
    @{ Row=27; C="11988776655"; D="" },
    @{ Row=28; C="11977665544"; D="" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "Vitor Ito"
    $ws.Cells.Item($row, 2).Value = 1578424633

    Set-TextValue $ws.Cells.Item($row, 3) $r.C
    if ($r.D -ne "") {
        Set-TextValue $ws.Cells.Item($row, 4) $r.D
    }

    $ws.Cells.Item($row, 5).Value = 1
    $ws.Cells.Item($row, 6).Value = 2
    $ws.Cells.Item($row, 7).Value = 3
    $ws.Cells.Item($row, 8).Value = 4
    $ws.Cells.Item($row, 9).Value = 5
    $ws.Cells.Item($row, 10).Value = 6
    $ws.Cells.Item($row, 11).Value = 7
    $ws.Cells.Item($row, 12).Value = 8
    $ws.Cells.Item($row, 13).Value = 9
    $ws.Cells.Item($row, 14).Value = 10
    $ws.Cells.Item($row, 15).Value = "Não"
}
